$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell without Excel coercing
# numeric-looking strings (e.g. "1.00", "0.999") into actual numbers.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '68.555.53'
$ws.Range("E2").Value = '  +0.69%  '

Set-TextValue $ws.Range("D3") '2.697.61'
$ws.Range("E3").Value = '  +2.18%  '

Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  -0.01%  '

Set-TextValue $ws.Range("D5") '599.13'
$ws.Range("E5").Value = '  +0.48%  '

Set-TextValue $ws.Range("D6") '160.16'
$ws.Range("E6").Value = '  +2.46%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  +0.23%  '

Set-TextValue $ws.Range("D9") '2.696.34'
$ws.Range("E9").Value = '  +2.15%  '

$ws.Range("E10").Value = '  -2.17%  '

$ws.Range("E11").Value = '  -0.37%  '

$ws.Range("E12").Value = '  +1.29%  '

$ws.Range("E13").Value = '  +2.56%  '

Set-TextValue $ws.Range("D14") '28.27'
$ws.Range("E14").Value = '  +1.00%  '

Set-TextValue $ws.Range("D15") '3.194.07'
$ws.Range("E15").Value = '  +2.23%  '

Set-TextValue $ws.Range("D16") '0.0000188'
$ws.Range("E16").Value = '  -0.60%  '

Set-TextValue $ws.Range("D17") '68.491.48'
$ws.Range("E17").Value = '  +0.62%  '

Set-TextValue $ws.Range("D18") '2.699.24'
$ws.Range("E18").Value = '  +2.20%  '

Set-TextValue $ws.Range("D19") '11.97'
$ws.Range("E19").Value = '  +5.62%  '

Set-TextValue $ws.Range("D20") '7.69'
$ws.Range("E20").Value = '  +4.08%  '

Set-TextValue $ws.Range("D21") '366.69'
$ws.Range("E21").Value = '  +1.26%  '

$ws.Range("E22").Value = '  +3.13%  '

$ws.Range("E23").Value = '  +2.33%  '

$ws.Range("E24").Value = '  +3.46%  '

Set-TextValue $ws.Range("D25") '74.52'
$ws.Range("E25").Value = '  -0.63%  '

Set-TextValue $ws.Range("D27") '10.10'
$ws.Range("E27").Value = '  +4.43%  '

Set-TextValue $ws.Range("D28") '2.834.09'
$ws.Range("E28").Value = '  +2.00%  '

$ws.Range("E29").Value = '  +1.01%  '

Set-TextValue $ws.Range("D30") '0.999'
$ws.Range("E30").Value = '  -0.24%  '

Set-TextValue $ws.Range("D31") '572.40'
$ws.Range("E31").Value = '  +3.50%  '

$ws.Range("E32").Value = '  +4.68%  '

Set-TextValue $ws.Range("D33") '8.24'
$ws.Range("E33").Value = '  +3.22%  '

$ws.Range("E34").Value = '  +5.67%  '

$ws.Range("E35").Value = '  +3.54%  '

$ws.Range("E36").Value = '  +6.79%  '

Set-TextValue $ws.Range("D37") '0.999'
$ws.Range("E37").Value = '  -0.05%  '

Set-TextValue $ws.Range("D38") '19.96'
$ws.Range("E38").Value = '  +3.06%  '

Set-TextValue $ws.Range("D39") '160.95'
$ws.Range("E39").Value = '  -0.11%  '

$ws.Range("E40").Value = '  +2.46%  '

$ws.Range("E41").Value = '  +2.08%  '

$ws.Range("E42").Value = '  +2.22%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D43") '2.66'
$ws.Range("E43").Value = '  +2.27%  '

$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range("D44") '17.87'
$ws.Range("E44").Value = '  +0.49%  '

$ws.Range("E45").Value = '  +0.05%  '

Set-TextValue $ws.Range("D46") '0.0₆0316'
$ws.Range("E46").Value = '  -6.60%  '

Set-TextValue $ws.Range("D47") '158.22'
$ws.Range("E47").Value = '  -0.39%  '

Set-TextValue $ws.Range("D48") '3.98'
$ws.Range("E48").Value = '  +6.83%  '

$ws.Range("E49").Value = '  +5.44%  '

Set-TextValue $ws.Range("D50") '0.601'
$ws.Range("E50").Value = '  +7.17%  '

Set-TextValue $ws.Range("D51") '22.08'
$ws.Range("E51").Value = '  +0.35%  '
